# Apply the "Red cedar" -> "Bois thermo-chauffé" rename, including the
# associated image filenames, and widen the two filename columns so the
# longer text fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the matiere + file name cells for the EB12 project (row 14).
$ws.Range("N14").Value = "Bois thermo-chauffé"
$ws.Range("R14").Value = "EB12-BOISTHERMOCHAUFFE-CAROTTE.jpg"
$ws.Range("T14").Value = "EB12-BOISTHERMOCHAUFFE-BASEDEDONNEES.jpg"

# Widen column R (nom_fichier_page_projets) and column T (nom_fichier_page_bdd)
# to accommodate the longer file names (values chosen to land as close as
# possible to the target OOXML column widths of 66.4375 / 54.0703 once the
# host snaps the character width to its pixel grid).
$ws.Columns.Item(18).ColumnWidth = 65.7142857142857
$ws.Columns.Item(20).ColumnWidth = 53.2857142857143
